# Southern Ocean - fixed deployment dates
# Fixed deployment dates based on cruise reports and WHOI documentation
#
# Moorings sheet, row 2 (GS05MOAS-GL524 glider deployment):
#   - Anchor Launch Time (F2) corrected
#   - Recover Date (G2) filled in

$wb = $excel.ActiveWorkbook

$moorings = $wb.Worksheets.Item("Moorings")

# Correct the anchor launch time (was 0.70138888888888884)
$moorings.Range("F2").Value = 0.82638888888888884

# Fill in the recover date (was blank)
$moorings.Range("G2").Value = 42471

# Leave the editor on the Moorings sheet, with the recover-date
# column in view, matching where the corrections were made.
[void]$moorings.Select()
[void]$moorings.Range("F8").Select()
